$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Make "Repayment Schedule" the active/selected sheet (was "Transactions" before)
$ws.Activate()

# A new (blank) column was inserted before the old "Outstanding"/"Late" columns (N),
# pushing the old N->O, O->P, P->Q
$ws.Columns("N").Insert()

# Update the selection on the now-active sheet
$ws.Range("T9").Select() | Out-Null
